$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 448, shifting existing rows 448..550 down to 449..551
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new record
$ws.Cells.Item(448, 1).Value = 5
$ws.Cells.Item(448, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(448, 3).Value = "Maule"
$ws.Cells.Item(448, 4).Value = 45204
$ws.Cells.Item(448, 5).Value = 7
$ws.Cells.Item(448, 6).Value = 100112003
$ws.Cells.Item(448, 7).Value = "Ajo"
$ws.Cells.Item(448, 8).Value = "Chino"
$ws.Cells.Item(448, 9).Value = "Primera"
$ws.Cells.Item(448, 10).Value = 200
$ws.Cells.Item(448, 11).Value = 21000
$ws.Cells.Item(448, 12).Value = 21000
$ws.Cells.Item(448, 13).Value = 21000
$ws.Cells.Item(448, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(448, 15).Value = "China"
$ws.Cells.Item(448, 16).Value = 2100
$ws.Cells.Item(448, 17).Value = 10
$ws.Cells.Item(448, 18).Value = "Hortaliza"
